# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> was "Office Theme" (only ever linked from the
#                             notes master's relationship)
#   ppt/theme/theme2.xml  -> was "Integral"      (the theme actually driving
#                             the slide master / presentation design)
#
# The authored change swaps the two themes' content in place, so the deck's
# visible design (slide master, i.e. every slide) switches from the green
# "Integral" palette over to the default blue/grey "Office" palette (and vice
# versa for the notes-only theme part). The font scheme and effect/format
# scheme are identical between the two themes already, so the entire visual
# delta is the 12-slot DrawingML colour scheme
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# Re-point the slide master's theme colour scheme (PowerPoint's
# ThemeColorScheme, which backs ppt/theme/theme2.xml - the theme that is
# actually applied to the slides) from the Integral palette to the Office
# palette so the presentation's design matches the target "Office Theme"
# colours.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Office Theme palette (RGB() long values = R + G*256 + B*65536), in the
# fixed ThemeColorScheme slot order 1-12:
#   1 dk1       000000
#   2 lt1       FFFFFF
#   3 dk2       44546A
#   4 lt2       E7E6E6
#   5 accent1   5B9BD5
#   6 accent2   ED7D31
#   7 accent3   A5A5A5
#   8 accent4   FFC000
#   9 accent5   4472C4
#  10 accent6   70AD47
#  11 hlink     0563C1
#  12 folHlink  954F72
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
